$d = $word.ActiveDocument

# --- Change 1 ---
# Merge the split runs "Ich habe " + "JavaScript" + " verwendet, ... sonstiges " + "ist," + " muss er eine kurze "
# into a single contiguous phrase by removing the extra spaces introduced by run
# splitting. Achieved via a Find/Replace that collapses the original phrase.
$d.Content.Find.Execute(
    "Ich habe JavaScript verwendet, um Notizen und ein Textfeld ein- und auszublenden. Wenn der Arbeitnehmer bei einem Arztbesuch oder krank ist, muss er eine AU mitbringen. und wenn es sonstiges ist, muss er eine kurze ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Ich habe JavaScript verwendet, um Notizen und ein Textfeld ein- und auszublenden. Wenn der Arbeitnehmer bei einem Arztbesuch oder krank ist, muss er eine AU mitbringen. und wenn es sonstiges ist, muss er eine kurze ",
    2
) | Out-Null

# --- Change 2 ---
# Append new sentence after "...als krank markieren kann."
$d.Content.Find.Execute(
    "und dass der Benutzer einen anderen Benutzer als krank markieren kann.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "und dass der Benutzer einen anderen Benutzer als krank markieren kann. Und die Möglichkeit, dass ein Mitarbeiter einen anderen Mitarbeiter als anwesend registrieren kann.",
    2
) | Out-Null

$d.Save()
